$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells for the season record columns (AD/AE/AF).
# Copy the existing header style (bold, bordered, centered) from A1
# onto the new header cells before writing their labels.
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Apply the team's season record to every player row (2 through 42)
for ($row = 2; $row -le 42; $row++) {
    $ws.Cells.Item($row, 30).Value = 93  # AD -> Wins
    $ws.Cells.Item($row, 31).Value = 69  # AE -> Losses
    $ws.Cells.Item($row, 32).Value = 0   # AF -> Ties
}
